# "Implementacao do brightness e do contraste"
# Update progress (%) tracking for alineas 1g, 1h and 1i on the Progress sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1g (row 15): progress revised down from 50% to 25%
$ws.Range("D15").Value = 25

# 1h (row 16, "brightness"): finished -> 100%; also clear the underline
# formatting that used to flag this row while it was still in review
$ws.Range("D16").Font.Underline = $false
$ws.Range("D16").Value = 100

# 1i (row 17, "contraste"): finished -> 100%
$ws.Range("D17").Value = 100

# Leave the selection where the user was last working
$ws.Range("D16").Select()
